$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.152424454689026
$ws.Range("B1").Value = 1.6400306224823
$ws.Range("C1").Value = 4.64898681640625
$ws.Range("D1").Value = 0.5950474739074707
$ws.Range("E1").Value = 0.6604103446006775
